# Update "想去人数" (column F) figures across the three affected sheets
# (展览, 演出, 全部类型) to match the latest data pull.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 5282
$ws1.Range("F6").Value  = 5282
$ws1.Range("F7").Value  = 158
$ws1.Range("F9").Value  = 533
$ws1.Range("F11").Value = 1187
$ws1.Range("F13").Value = 5218
$ws1.Range("F14").Value = 30
$ws1.Range("F16").Value = 92
$ws1.Range("F17").Value = 2306
$ws1.Range("F18").Value = 2306
$ws1.Range("F19").Value = 258
$ws1.Range("F22").Value = 3922
$ws1.Range("F26").Value = 3850
$ws1.Range("F29").Value = 250
$ws1.Range("F32").Value = 113
$ws1.Range("F35").Value = 139
$ws1.Range("F36").Value = 23
$ws1.Range("F37").Value = 6840
$ws1.Range("F38").Value = 1111
$ws1.Range("F39").Value = 530
$ws1.Range("F41").Value = 62
$ws1.Range("F42").Value = 1393
$ws1.Range("F43").Value = 175
$ws1.Range("F44").Value = 701
$ws1.Range("F46").Value = 2330
$ws1.Range("F47").Value = 316
$ws1.Range("F49").Value = 8
$ws1.Range("F50").Value = 787

# --- Sheet: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value  = 5
$ws2.Range("F11").Value = 66
$ws2.Range("F14").Value = 6
$ws2.Range("F16").Value = 9

# --- Sheet: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 5
$ws4.Range("F7").Value  = 5283
$ws4.Range("F8").Value  = 5283
$ws4.Range("F9").Value  = 158
$ws4.Range("F12").Value = 533
$ws4.Range("F14").Value = 1187
$ws4.Range("F16").Value = 30
$ws4.Range("F18").Value = 92
$ws4.Range("F19").Value = 2306
$ws4.Range("F20").Value = 2307
$ws4.Range("F21").Value = 258
$ws4.Range("F24").Value = 3922
$ws4.Range("F25").Value = 3850
$ws4.Range("F28").Value = 250
$ws4.Range("F31").Value = 113
$ws4.Range("F33").Value = 139
$ws4.Range("F34").Value = 23
$ws4.Range("F36").Value = 6840
$ws4.Range("F37").Value = 1111
$ws4.Range("F38").Value = 530
$ws4.Range("F41").Value = 62
$ws4.Range("F42").Value = 1393
$ws4.Range("F43").Value = 175
$ws4.Range("F44").Value = 701
$ws4.Range("F46").Value = 2330
$ws4.Range("F47").Value = 316
$ws4.Range("F49").Value = 787
